$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 6912735.5   # H116
$ws.Cells.Item(116, 9).Value = 8335422   # I116
$ws.Cells.Item(116, 11).Value = 8335422   # K116
$ws.Cells.Item(116, 13).Value = -8331980   # M116
$ws.Cells.Item(125, 8).Value = 1430.7   # H125
$ws.Cells.Item(125, 9).Value = 633.3333   # I125
$ws.Cells.Item(125, 10).Value = 1772.4286   # J125
$ws.Cells.Item(125, 11).Value = 5699.9997   # K125
$ws.Cells.Item(125, 12).Value = 15951.8574   # L125
$ws.Cells.Item(125, 13).Value = -3239.9997   # M125
$ws.Cells.Item(125, 14).Value = -20871.8574   # N125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8943.464   # H32
$ws.Cells.Item(32, 9).Value = 2614.8965   # I32
$ws.Cells.Item(32, 10).Value = 24237.5   # J32
$ws.Cells.Item(32, 11).Value = 2614.8965   # K32
$ws.Cells.Item(32, 12).Value = 24237.5   # L32
$ws.Cells.Item(32, 13).Value = -2327.8965   # M32
$ws.Cells.Item(32, 14).Value = -24811.5   # N32
$ws.Cells.Item(102, 8).Value = 2444   # H102
$ws.Cells.Item(102, 9).Value = 2444   # I102
$ws.Cells.Item(102, 11).Value = 2444   # K102
$ws.Cells.Item(102, 13).Value = -822   # M102

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1666.6666   # H99
$ws.Cells.Item(99, 9).Value = 1666.6666   # I99
$ws.Cells.Item(99, 11).Value = 1666.6666   # K99
$ws.Cells.Item(99, 13).Value = -168.6666   # M99
$ws.Cells.Item(134, 8).Value = 3676.8235   # H134
$ws.Cells.Item(134, 9).Value = 2433.8   # I134
$ws.Cells.Item(134, 11).Value = 7301.400000000001   # K134
$ws.Cells.Item(134, 13).Value = -4766.400000000001   # M134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10555109   # H31
$ws.Cells.Item(31, 9).Value = 18273752   # I31
$ws.Cells.Item(31, 10).Value = 7410477   # J31
$ws.Cells.Item(31, 11).Value = 18273752   # K31
$ws.Cells.Item(31, 12).Value = 7410477   # L31
$ws.Cells.Item(31, 13).Value = -18273457   # M31
$ws.Cells.Item(31, 14).Value = -7411067   # N31
$ws.Cells.Item(34, 8).Value = 10555109   # H34
$ws.Cells.Item(34, 9).Value = 18273752   # I34
$ws.Cells.Item(34, 10).Value = 7410477   # J34
$ws.Cells.Item(34, 11).Value = 18273752   # K34
$ws.Cells.Item(34, 12).Value = 7410477   # L34
$ws.Cells.Item(34, 13).Value = -18273550   # M34
$ws.Cells.Item(34, 14).Value = -7410881   # N34
$ws.Cells.Item(99, 8).Value = 3579990.8   # H99
$ws.Cells.Item(99, 9).Value = 7150501.5   # I99
$ws.Cells.Item(99, 10).Value = 9480   # J99
$ws.Cells.Item(99, 11).Value = 7150501.5   # K99
$ws.Cells.Item(99, 12).Value = 9480   # L99
$ws.Cells.Item(99, 13).Value = -7149003.5   # M99
$ws.Cells.Item(99, 14).Value = -12476   # N99
$ws.Cells.Item(126, 8).Value = 3579990.8   # H126
$ws.Cells.Item(126, 9).Value = 7150501.5   # I126
$ws.Cells.Item(126, 10).Value = 9480   # J126
$ws.Cells.Item(126, 11).Value = 21451504.5   # K126
$ws.Cells.Item(126, 12).Value = 28440   # L126
$ws.Cells.Item(126, 13).Value = -21449034.5   # M126
$ws.Cells.Item(126, 14).Value = -33380   # N126

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 2000   # H18
$ws.Cells.Item(18, 10).Value = 2000   # J18
$ws.Cells.Item(18, 12).Value = 2000   # L18
$ws.Cells.Item(18, 14).Value = -2586   # N18
$ws.Cells.Item(122, 8).Value = 2356.4375   # H122
$ws.Cells.Item(122, 9).Value = 2491.9167   # I122
$ws.Cells.Item(122, 10).Value = 1950   # J122
$ws.Cells.Item(122, 11).Value = 7475.750100000001   # K122
$ws.Cells.Item(122, 12).Value = 5850   # L122
$ws.Cells.Item(122, 13).Value = -5025.750100000001   # M122
$ws.Cells.Item(122, 14).Value = -10750   # N122
$ws.Cells.Item(132, 8).Value = 2577.0588   # H132
$ws.Cells.Item(132, 9).Value = 2037.4546   # I132
$ws.Cells.Item(132, 10).Value = 3566.3333   # J132
$ws.Cells.Item(132, 11).Value = 6112.3638   # K132
$ws.Cells.Item(132, 12).Value = 10698.9999   # L132
$ws.Cells.Item(132, 13).Value = -3582.3638   # M132
$ws.Cells.Item(132, 14).Value = -15758.9999   # N132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 166.85   # H55
$ws.Cells.Item(55, 10).Value = 235.71428   # J55
$ws.Cells.Item(55, 12).Value = 235.71428   # L55
$ws.Cells.Item(55, 14).Value = -581.71428   # N55
$ws.Cells.Item(61, 8).Value = 1187.875   # H61
$ws.Cells.Item(61, 9).Value = 899.6   # I61
$ws.Cells.Item(61, 10).Value = 1668.3334   # J61
$ws.Cells.Item(61, 11).Value = 899.6   # K61
$ws.Cells.Item(61, 12).Value = 1668.3334   # L61
$ws.Cells.Item(61, 13).Value = -697.6   # M61
$ws.Cells.Item(61, 14).Value = -2072.3334   # N61
$ws.Cells.Item(82, 8).Value = 2196.389   # H82
$ws.Cells.Item(82, 9).Value = 1512.375   # I82
$ws.Cells.Item(82, 10).Value = 2743.6   # J82
$ws.Cells.Item(82, 11).Value = 1512.375   # K82
$ws.Cells.Item(82, 12).Value = 2743.6   # L82
$ws.Cells.Item(82, 13).Value = -1151.375   # M82
$ws.Cells.Item(82, 14).Value = -3465.6   # N82
$ws.Cells.Item(85, 8).Value = 2196.389   # H85
$ws.Cells.Item(85, 9).Value = 1512.375   # I85
$ws.Cells.Item(85, 10).Value = 2743.6   # J85
$ws.Cells.Item(85, 11).Value = 1512.375   # K85
$ws.Cells.Item(85, 12).Value = 2743.6   # L85
$ws.Cells.Item(85, 13).Value = -264.375   # M85
$ws.Cells.Item(85, 14).Value = -5239.6   # N85
$ws.Cells.Item(93, 8).Value = 1701.2106   # H93
$ws.Cells.Item(93, 9).Value = 1594.9375   # I93
$ws.Cells.Item(93, 10).Value = 2268   # J93
$ws.Cells.Item(93, 11).Value = 1594.9375   # K93
$ws.Cells.Item(93, 12).Value = 2268   # L93
$ws.Cells.Item(93, 13).Value = -346.9375   # M93
$ws.Cells.Item(93, 14).Value = -4764   # N93
$ws.Cells.Item(113, 8).Value = 1187.875   # H113
$ws.Cells.Item(113, 9).Value = 899.6   # I113
$ws.Cells.Item(113, 10).Value = 1668.3334   # J113
$ws.Cells.Item(113, 11).Value = 899.6   # K113
$ws.Cells.Item(113, 12).Value = 1668.3334   # L113
$ws.Cells.Item(113, 13).Value = 1270.4   # M113
$ws.Cells.Item(113, 14).Value = -6008.3334   # N113
$ws.Cells.Item(132, 8).Value = 3692929   # H132
$ws.Cells.Item(132, 9).Value = 4482421   # I132
$ws.Cells.Item(132, 10).Value = 8633.333000000001   # J132
$ws.Cells.Item(132, 11).Value = 13447263   # K132
$ws.Cells.Item(132, 12).Value = 25899.999   # L132
$ws.Cells.Item(132, 13).Value = -13444733   # M132
$ws.Cells.Item(132, 14).Value = -30959.999   # N132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 31158   # H21
$ws.Cells.Item(21, 9).Value = 2013.75   # I21
$ws.Cells.Item(21, 11).Value = 2013.75   # K21
$ws.Cells.Item(21, 13).Value = -1778.75   # M21
$ws.Cells.Item(35, 8).Value = 31158   # H35
$ws.Cells.Item(35, 9).Value = 2013.75   # I35
$ws.Cells.Item(35, 11).Value = 2013.75   # K35
$ws.Cells.Item(35, 13).Value = -1723.75   # M35
$ws.Cells.Item(62, 8).Value = 6695.3335   # H62
$ws.Cells.Item(62, 9).Value = 4651   # I62
$ws.Cells.Item(62, 10).Value = 7279.4287   # J62
$ws.Cells.Item(62, 11).Value = 4651   # K62
$ws.Cells.Item(62, 12).Value = 7279.4287   # L62
$ws.Cells.Item(62, 13).Value = -4027   # M62
$ws.Cells.Item(62, 14).Value = -8527.4287   # N62
$ws.Cells.Item(65, 8).Value = 6695.3335   # H65
$ws.Cells.Item(65, 9).Value = 4651   # I65
$ws.Cells.Item(65, 10).Value = 7279.4287   # J65
$ws.Cells.Item(65, 11).Value = 23255   # K65
$ws.Cells.Item(65, 12).Value = 36397.14350000001   # L65
$ws.Cells.Item(65, 13).Value = -20135   # M65
$ws.Cells.Item(65, 14).Value = -42637.14350000001   # N65
$ws.Cells.Item(81, 8).Value = 71430800   # H81
$ws.Cells.Item(81, 10).Value = 2485.25   # J81
$ws.Cells.Item(81, 12).Value = 4970.5   # L81
$ws.Cells.Item(81, 14).Value = -7092.5   # N81
$ws.Cells.Item(84, 8).Value = 71430800   # H84
$ws.Cells.Item(84, 10).Value = 2485.25   # J84
$ws.Cells.Item(84, 12).Value = 24852.5   # L84
$ws.Cells.Item(84, 14).Value = -35460.5   # N84
$ws.Cells.Item(122, 8).Value = 52634064   # H122
$ws.Cells.Item(122, 9).Value = 90911010   # I122
$ws.Cells.Item(122, 10).Value = 3275   # J122
$ws.Cells.Item(122, 11).Value = 272733030   # K122
$ws.Cells.Item(122, 12).Value = 9825   # L122
$ws.Cells.Item(122, 13).Value = -272730580   # M122
$ws.Cells.Item(122, 14).Value = -14725   # N122
$ws.Cells.Item(126, 8).Value = 4830.4   # H126
$ws.Cells.Item(126, 9).Value = 6472   # I126
$ws.Cells.Item(126, 10).Value = 1000   # J126
$ws.Cells.Item(126, 11).Value = 19416   # K126
$ws.Cells.Item(126, 12).Value = 3000   # L126
$ws.Cells.Item(126, 13).Value = -16946   # M126
$ws.Cells.Item(126, 14).Value = -7940   # N126
